$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.143.82"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.430.18"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'318.18"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'89.62"
$ws.Range("E6").Value = "  -2.80%  "
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D11").Value = "'32.10"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "2.803.86"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").Value = "'15.67"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").Value = "2.424.94"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "'0.778"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("D18").Value = "41.081.14"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").Value = "0.0₃0928"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("E20").Value = "  -2.43%  "
$ws.Range("D21").Value = "'72.37"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "'11.08"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("D23").Value = "'235.11"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "'1.87"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").Value = "'24.20"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").Value = "'9.64"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").Value = "'34.61"
$ws.Range("E30").Value = "  -3.96%  "
$ws.Range("D31").Value = "'158.92"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").Value = "'5.27"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "'17.07"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.97"
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("D42").Value = "'2.32"
$ws.Range("E42").Value = "  -6.06%  "
$ws.Range("D43").Value = "1.996.21"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").Value = "'18.58"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("D47").Value = "'9.57"
$ws.Range("E47").Value = "  +4.28%  "
$ws.Range("D48").Value = "2.664.94"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "'94.90"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").Value = "'51.86"
$ws.Range("E51").Value = "  -0.81%  "
